# Actualización automática de tasas-transfi.xlsx
# - Actualiza el resumen de conversión del día en la hoja "Hoja1" (A1)
# - Actualiza las tasas (N10/O10/N12/O12) en la hoja "tasas"

$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: texto de conversión del día ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$conversionText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.88 = 6868.23 pesos`n✅ 6868.23 pesos = 1.87 = 916.76 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $conversionText

# --- tasas: tasas de cambio actualizadas ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 530.85
$ws2.Range("O10").Value = 3646
$ws2.Range("N12").Value = 3670.99
$ws2.Range("O12").Value = 490
